$d = $word.ActiveDocument

# --- Locate the " - organizacija 2 " chunk inside the title paragraph ---
$findRng = $d.Content
$null = $findRng.Find.Execute(" – organizacija 2 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$run1Start = $findRng.Start

# --- Step 1: shrink that chunk down to a single space " " ---
$findRng.Text = " "

# Position right after the remaining single space.
$afterSpace = $run1Start + 1

# --- Step 2: insert the capitalised, replacement text right after that space ---
$newChunk = [char]8211 + " Organizacija 2 "
$insPoint = $d.Range($afterSpace, $afterSpace)
$insPoint.InsertAfter($newChunk)

# --- Step 3: italicise (and explicitly un-bold) everything from the newly inserted
#     text through to the end of the paragraph text (covers the trailing "- ponavljanje") ---
$para = $d.Paragraphs.Item(1)
$paraTextEnd = $para.Range.End - 1
$italicRange = $d.Range($afterSpace, $paraTextEnd)
$italicRange.Font.Italic = $true
$italicRange.Font.Bold = $false

# --- Step 4: move the "_GoBack" bookmark so that it starts right after the space
#     and spans through to the end of the paragraph text ---
$bmRange = $d.Range($afterSpace, $paraTextEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
